$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set the "Name" value (B4) to "ProfessionVs"
$ws.Range("B4").Value = "ProfessionVs"

# Update the "Date" value (B8) to the new timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
